$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-6: update financial figures (per-cell) ---
$ws.Range("D2").Value = 1821
$ws.Range("E2").Value = 359
$ws.Range("F2").Value = 359
$ws.Range("G2").Value = 333
$ws.Range("H2").Value = 252
$ws.Range("I2").Value = 252
$ws.Range("K2").Value = 5126
$ws.Range("L2").Value = 2016
$ws.Range("M2").Value = 3110
$ws.Range("N2").Value = 3110
$ws.Range("P2").Value = 410
$ws.Range("Q2").Value = 341
$ws.Range("R2").Value = 21
$ws.Range("S2").Value = -226
$ws.Range("T2").Value = 34
$ws.Range("U2").Value = 307
$ws.Range("V2").Value = 975
$ws.Range("W2").Value = 19.7
$ws.Range("X2").Value = 13.82
$ws.Range("Y2").Value = 8.42
$ws.Range("Z2").Value = 4.92
$ws.Range("AA2").Value = 64.81
$ws.Range("AB2").Value = 673.22
$ws.Range("AC2").Value = 614
$ws.Range("AE2").Value = 8236
$ws.Range("AF2").Value = 0.73
$ws.Range("AG2").Value = 30
$ws.Range("AH2").Value = 0.5
$ws.Range("AI2").Value = 4.5
$ws.Range("AJ2").Value = 41040895
$ws.Range("D3").Value = 1870
$ws.Range("E3").Value = 355
$ws.Range("F3").Value = 355
$ws.Range("G3").Value = 338
$ws.Range("H3").Value = 257
$ws.Range("I3").Value = 257
$ws.Range("K3").Value = 5279
$ws.Range("L3").Value = 1915
$ws.Range("M3").Value = 3364
$ws.Range("N3").Value = 3364
$ws.Range("P3").Value = 410
$ws.Range("Q3").Value = 286
$ws.Range("R3").Value = -89
$ws.Range("S3").Value = -69
$ws.Range("T3").Value = 24
$ws.Range("U3").Value = 262
$ws.Range("V3").Value = 917
$ws.Range("W3").Value = 18.96
$ws.Range("X3").Value = 13.76
$ws.Range("Y3").Value = 7.95
$ws.Range("Z3").Value = 4.95
$ws.Range("AA3").Value = 56.94
$ws.Range("AB3").Value = 735
$ws.Range("AC3").Value = 627
$ws.Range("AD3").Value = 9.199999999999999
$ws.Range("AE3").Value = 8907
$ws.Range("AF3").Value = 0.65
$ws.Range("AG3").Value = 30
$ws.Range("AH3").Value = 0.52
$ws.Range("AI3").Value = 4.4
$ws.Range("AJ3").Value = 41040895
$ws.Range("D4").Value = 1830
$ws.Range("E4").Value = 369
$ws.Range("F4").Value = 369
$ws.Range("G4").Value = 363
$ws.Range("H4").Value = 271
$ws.Range("I4").Value = 271
$ws.Range("K4").Value = 5534
$ws.Range("L4").Value = 1909
$ws.Range("M4").Value = 3625
$ws.Range("N4").Value = 3625
$ws.Range("P4").Value = 410
$ws.Range("Q4").Value = 493
$ws.Range("R4").Value = -77
$ws.Range("S4").Value = -198
$ws.Range("T4").Value = 38
$ws.Range("U4").Value = 454
$ws.Range("V4").Value = 731
$ws.Range("W4").Value = 20.16
$ws.Range("X4").Value = 14.8
$ws.Range("Y4").Value = 7.75
$ws.Range("Z4").Value = 5.01
$ws.Range("AA4").Value = 52.66
$ws.Range("AB4").Value = 798.5700000000001
$ws.Range("AC4").Value = 660
$ws.Range("AD4").Value = 7.97
$ws.Range("AE4").Value = 9598
$ws.Range("AF4").Value = 0.55
$ws.Range("AG4").Value = 30
$ws.Range("AH4").Value = 0.57
$ws.Range("AI4").Value = 4.18
$ws.Range("AJ4").Value = 41040895
$ws.Range("D5").Value = 1784
$ws.Range("E5").Value = 362
$ws.Range("F5").Value = 362
$ws.Range("G5").Value = 382
$ws.Range("H5").Value = 289
$ws.Range("I5").Value = 289
$ws.Range("K5").Value = 5626
$ws.Range("L5").Value = 1728
$ws.Range("M5").Value = 3898
$ws.Range("N5").Value = 3898
$ws.Range("P5").Value = 410
$ws.Range("Q5").Value = 314
$ws.Range("R5").Value = -89
$ws.Range("S5").Value = -164
$ws.Range("T5").Value = 14
$ws.Range("U5").Value = 300
$ws.Range("V5").Value = 578
$ws.Range("W5").Value = 20.29
$ws.Range("X5").Value = 16.17
$ws.Range("Y5").Value = 7.67
$ws.Range("Z5").Value = 5.17
$ws.Range("AA5").Value = 44.32
$ws.Range("AB5").Value = 865.22
$ws.Range("AC5").Value = 703
$ws.Range("AD5").Value = 7
$ws.Range("AE5").Value = 10322
$ws.Range("AF5").Value = 0.48
$ws.Range("AG5").Value = 50
$ws.Range("AH5").Value = 1.02
$ws.Range("AI5").Value = 6.54
$ws.Range("AJ5").Value = 41040895
$ws.Range("D6").Value = 1740
$ws.Range("E6").Value = 339
$ws.Range("F6").Value = 339
$ws.Range("G6").Value = 346
$ws.Range("H6").Value = 254
$ws.Range("I6").Value = 254
$ws.Range("K6").Value = 5758
$ws.Range("L6").Value = 1412
$ws.Range("M6").Value = 4346
$ws.Range("N6").Value = 4346
$ws.Range("P6").Value = 410
$ws.Range("Q6").Value = 93
$ws.Range("R6").Value = -210
$ws.Range("S6").Value = -154
$ws.Range("T6").Value = 279
$ws.Range("U6").Value = -186
$ws.Range("V6").Value = 443
$ws.Range("W6").Value = 19.5
$ws.Range("X6").Value = 14.62
$ws.Range("Y6").Value = 6.17
$ws.Range("Z6").Value = 4.47
$ws.Range("AA6").Value = 32.48
$ws.Range("AB6").Value = 922.21
$ws.Range("AC6").Value = 619
$ws.Range("AD6").Value = 5.55
$ws.Range("AE6").Value = 11509
$ws.Range("AF6").Value = 0.3
$ws.Range("AG6").Value = 50
$ws.Range("AH6").Value = 1.46
$ws.Range("AI6").Value = 7.43
$ws.Range("AJ6").Value = 41040895

# --- Rows 2-5: clear columns J (당기순이익(비지배)) and O (자본총계(비지배)) ---
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# --- Rows 7-9: clear all financial data except index/label columns (A, B, C) ---
$ws.Range("D7:AJ9").ClearContents()
